# Auto-generated Excel COM-interop edit script
# Updates FFXIV Leve profit-calculation snapshot cells (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across all eight market sheets to reflect
# the latest scheduled market-data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3498.7856
$ws.Range("I40").Value = 2881.9412
$ws.Range("J40").Value = 4452.091
$ws.Range("K40").Value = 2881.9412
$ws.Range("L40").Value = 4452.091
$ws.Range("M40").Value = -2706.9412
$ws.Range("N40").Value = -4802.091
$ws.Range("H98").Value = 790.8823
$ws.Range("I98").Value = 629.73334
$ws.Range("K98").Value = 629.73334
$ws.Range("M98").Value = 868.26666
$ws.Range("H100").Value = 1465.8
$ws.Range("I100").Value = 1349.3334
$ws.Range("K100").Value = 1349.3334
$ws.Range("M100").Value = -808.3334
$ws.Range("H122").Value = 790.8823
$ws.Range("I122").Value = 629.73334
$ws.Range("K122").Value = 1889.20002
$ws.Range("M122").Value = 560.79998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3584.3333
$ws.Range("I2").Value = 2801.2
$ws.Range("K2").Value = 2801.2
$ws.Range("M2").Value = -2688.2
$ws.Range("H75").Value = 80000
$ws.Range("I75").Value = 80000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 80000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -79126
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 80000
$ws.Range("I78").Value = 80000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 240000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -235632
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 65025
$ws.Range("J80").Value = 66666.664
$ws.Range("L80").Value = 66666.664
$ws.Range("N80").Value = -68662.664
$ws.Range("H83").Value = 65025
$ws.Range("J83").Value = 66666.664
$ws.Range("L83").Value = 199999.992
$ws.Range("N83").Value = -209983.992
$ws.Range("H86").Value = 99994.5
$ws.Range("I86").Value = 99994.5
$ws.Range("K86").Value = 99994.5
$ws.Range("M86").Value = -98808.5
$ws.Range("H89").Value = 99994.5
$ws.Range("I89").Value = 99994.5
$ws.Range("K89").Value = 299983.5
$ws.Range("M89").Value = -294055.5
$ws.Range("H116").Value = 3584.3333
$ws.Range("I116").Value = 2801.2
$ws.Range("K116").Value = 2801.2
$ws.Range("M116").Value = -507.1999999999998
$ws.Range("H122").Value = 2377.5386
$ws.Range("J122").Value = 2767.5715
$ws.Range("L122").Value = 8302.7145
$ws.Range("N122").Value = -13202.7145
$ws.Range("H132").Value = 36501.87
$ws.Range("I132").Value = 2518.9636
$ws.Range("J132").Value = 161105.86
$ws.Range("K132").Value = 7556.8908
$ws.Range("L132").Value = 483317.58
$ws.Range("M132").Value = -5026.8908
$ws.Range("N132").Value = -488377.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3584.3333
$ws.Range("I3").Value = 2801.2
$ws.Range("K3").Value = 2801.2
$ws.Range("M3").Value = -2687.2
$ws.Range("H20").Value = 1887.3529
$ws.Range("I20").Value = 1755.75
$ws.Range("K20").Value = 1755.75
$ws.Range("M20").Value = -1508.75
$ws.Range("H80").Value = 603.9
$ws.Range("I80").Value = 540
$ws.Range("J80").Value = 667.8
$ws.Range("K80").Value = 540
$ws.Range("L80").Value = 667.8
$ws.Range("M80").Value = 458
$ws.Range("N80").Value = -2663.8
$ws.Range("H83").Value = 603.9
$ws.Range("I83").Value = 540
$ws.Range("J83").Value = 667.8
$ws.Range("K83").Value = 2700
$ws.Range("L83").Value = 3339
$ws.Range("M83").Value = 2292
$ws.Range("N83").Value = -13323
$ws.Range("H107").Value = 3127.1052
$ws.Range("J107").Value = 2953.6
$ws.Range("L107").Value = 2953.6
$ws.Range("N107").Value = -6793.6
$ws.Range("H134").Value = 1725.5
$ws.Range("I134").Value = 1667.8
$ws.Range("K134").Value = 5003.4
$ws.Range("M134").Value = -2468.4
$ws.Range("H141").Value = 54492
$ws.Range("I141").Value = 19000
$ws.Range("J141").Value = 89984
$ws.Range("K141").Value = 19000
$ws.Range("L141").Value = 89984
$ws.Range("M141").Value = -13820
$ws.Range("N141").Value = -100344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11460952
$ws.Range("I58").Value = 2239.5
$ws.Range("J58").Value = 19645746
$ws.Range("K58").Value = 2239.5
$ws.Range("L58").Value = 19645746
$ws.Range("M58").Value = -2036.5
$ws.Range("N58").Value = -19646152
$ws.Range("H86").Value = 13749.75
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 13749.75
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws.Range("H105").Value = 6791.8125
$ws.Range("I105").Value = 6678.6
$ws.Range("K105").Value = 6678.6
$ws.Range("M105").Value = -4931.6
$ws.Range("H131").Value = 31517.572
$ws.Range("J131").Value = 31517.572
$ws.Range("L131").Value = 31517.572
$ws.Range("N131").Value = -41597.572
$ws.Range("H132").Value = 3931.1052
$ws.Range("I132").Value = 3587.2964
$ws.Range("J132").Value = 4775
$ws.Range("K132").Value = 10761.8892
$ws.Range("L132").Value = 14325
$ws.Range("M132").Value = -8231.889200000001
$ws.Range("N132").Value = -19385
$ws.Range("H134").Value = 2979.4167
$ws.Range("I134").Value = 2507
$ws.Range("J134").Value = 3924.25
$ws.Range("K134").Value = 7521
$ws.Range("L134").Value = 11772.75
$ws.Range("M134").Value = -4986
$ws.Range("N134").Value = -16842.75
$ws.Range("H135").Value = 98365.71000000001
$ws.Range("J135").Value = 98365.71000000001
$ws.Range("L135").Value = 98365.71000000001
$ws.Range("N135").Value = -108505.71
$ws.Range("H136").Value = 11460952
$ws.Range("I136").Value = 2239.5
$ws.Range("J136").Value = 19645746
$ws.Range("K136").Value = 6718.5
$ws.Range("L136").Value = 58937238
$ws.Range("M136").Value = -4168.5
$ws.Range("N136").Value = -58942338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 930
$ws.Range("I12").Value = 1001
$ws.Range("J12").Value = 918.1667
$ws.Range("K12").Value = 3003
$ws.Range("L12").Value = 2754.5001
$ws.Range("M12").Value = -2830
$ws.Range("N12").Value = -3100.5001
$ws.Range("H33").Value = 160.26923
$ws.Range("I33").Value = 115.14286
$ws.Range("J33").Value = 212.91667
$ws.Range("K33").Value = 690.85716
$ws.Range("L33").Value = 1277.50002
$ws.Range("M33").Value = -407.85716
$ws.Range("N33").Value = -1843.50002
$ws.Range("H55").Value = 1002326.7
$ws.Range("J55").Value = 4997.5
$ws.Range("L55").Value = 14992.5
$ws.Range("N55").Value = -15346.5
$ws.Range("H80").Value = 6333
$ws.Range("I80").Value = 4999.5
$ws.Range("K80").Value = 14998.5
$ws.Range("M80").Value = -14062.5
$ws.Range("H83").Value = 6333
$ws.Range("I83").Value = 4999.5
$ws.Range("K83").Value = 44995.5
$ws.Range("M83").Value = -40315.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9835.333000000001
$ws.Range("I70").Value = 9835.333000000001
$ws.Range("K70").Value = 9835.333000000001
$ws.Range("M70").Value = -9565.333000000001
$ws.Range("H73").Value = 9835.333000000001
$ws.Range("I73").Value = 9835.333000000001
$ws.Range("K73").Value = 9835.333000000001
$ws.Range("M73").Value = -8899.333000000001
$ws.Range("H107").Value = 1128.2858
$ws.Range("I107").Value = 483.16666
$ws.Range("J107").Value = 4999
$ws.Range("K107").Value = 483.16666
$ws.Range("L107").Value = 4999
$ws.Range("M107").Value = 1436.83334
$ws.Range("N107").Value = -8839
$ws.Range("H132").Value = 1099.3334
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 1399
$ws.Range("K132").Value = 1500
$ws.Range("L132").Value = 4197
$ws.Range("M132").Value = 1030
$ws.Range("N132").Value = -9257

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2085.2
$ws.Range("I122").Value = 1959.5385
$ws.Range("J122").Value = 2902
$ws.Range("K122").Value = 5878.6155
$ws.Range("L122").Value = 8706
$ws.Range("M122").Value = -3428.6155
$ws.Range("N122").Value = -13606
$ws.Range("H132").Value = 4137.653
$ws.Range("I132").Value = 2277.1538
$ws.Range("K132").Value = 6831.4614
$ws.Range("M132").Value = -4301.4614
$ws.Range("H136").Value = 2593.9565
$ws.Range("J136").Value = 3750
$ws.Range("L136").Value = 11250
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1272.8
$ws.Range("I132").Value = 1272.8
$ws.Range("K132").Value = 3818.4
$ws.Range("M132").Value = -1288.4
$ws.Range("H136").Value = 5392.727
$ws.Range("I136").Value = 3259.353
$ws.Range("K136").Value = 9778.059000000001
$ws.Range("M136").Value = -7228.059000000001

